$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in score data for row 2 (E. Wira)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 7

# Fill in score data for row 3 (F. Hibban)
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 7

# Fill in score data for row 8 (L. Cahya)
$ws.Range("C8").Value = 7

# Fill in score data for row 9 (M. Afif)
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 6

# Update the selected cell/active cell on the sheet
$ws.Range("F3").Select()

$wb.Save()
